# Update "想去人数" (want-to-go count) values in the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 775
    $ws.Range("F3").Value = 4187
    $ws.Range("F5").Value = 760
}
